$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Forecasted Consumption (MW) values for rows 2-97
$aValues = @(
    5510, 5460, 5430, 5400, 5360, 5350, 5350, 5340, 5330, 5330, 5330, 5340, 5350, 5360, 5370, 5400, 5450, 5520, 5620, 5740, 5870, 6040, 6210, 6370, 6640, 6810, 6950, 7070, 7150, 7180, 7180, 7170, 7140, 7070, 6960, 6840, 6700, 6550, 6410, 6240, 6060, 5930, 5820, 5730, 5650, 5600, 5580, 5570, 5580, 5590, 5600, 5620, 5660, 5700, 5740, 5780, 5830, 5870, 5930, 5990, 6070, 6170, 6280, 6400, 6520, 6630, 6740, 6840, 6960, 7050, 7170, 7300, 7460, 7570, 7630, 7620, 7560, 7510, 7450, 7400, 7320, 7220, 7130, 6970, 6730, 6580, 6420, 6280, 6140, 6000, 5890, 5780, 5650, 5580, 5550, 5500
)

# Timestamp (Excel serial date) values for rows 2-97
$bValues = @(
    45946, 45946.01041666666, 45946.02083333334, 45946.03125, 45946.04166666666, 45946.05208333334, 45946.0625, 45946.07291666666, 45946.08333333334, 45946.09375, 45946.10416666666, 45946.11458333334, 45946.125, 45946.13541666666, 45946.14583333334, 45946.15625, 45946.16666666666, 45946.17708333334, 45946.1875, 45946.19791666666, 45946.20833333334, 45946.21875, 45946.22916666666, 45946.23958333334, 45946.25, 45946.26041666666, 45946.27083333334, 45946.28125, 45946.29166666666, 45946.30208333334, 45946.3125, 45946.32291666666, 45946.33333333334, 45946.34375, 45946.35416666666, 45946.36458333334, 45946.375, 45946.38541666666, 45946.39583333334, 45946.40625, 45946.41666666666, 45946.42708333334, 45946.4375, 45946.44791666666, 45946.45833333334, 45946.46875, 45946.47916666666, 45946.48958333334, 45946.5, 45946.51041666666, 45946.52083333334, 45946.53125, 45946.54166666666, 45946.55208333334, 45946.5625, 45946.57291666666, 45946.58333333334, 45946.59375, 45946.60416666666, 45946.61458333334, 45946.625, 45946.63541666666, 45946.64583333334, 45946.65625, 45946.66666666666, 45946.67708333334, 45946.6875, 45946.69791666666, 45946.70833333334, 45946.71875, 45946.72916666666, 45946.73958333334, 45946.75, 45946.76041666666, 45946.77083333334, 45946.78125, 45946.79166666666, 45946.80208333334, 45946.8125, 45946.82291666666, 45946.83333333334, 45946.84375, 45946.85416666666, 45946.86458333334, 45946.875, 45946.88541666666, 45946.89583333334, 45946.90625, 45946.91666666666, 45946.92708333334, 45946.9375, 45946.94791666666, 45946.95833333334, 45946.96875, 45946.97916666666, 45946.98958333334
)

$rowCount = $aValues.Length

# Build 2D arrays for bulk write into the worksheet
$aArr = New-Object "object[,]" $rowCount,1
$bArr = New-Object "object[,]" $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $aArr[$i,0] = $aValues[$i]
    $bArr[$i,0] = $bValues[$i]
}

$startRow = 2
$endRow = $startRow + $rowCount - 1

$aRange = $ws.Range("A$startRow`:A$endRow")
$bRange = $ws.Range("B$startRow`:B$endRow")

$aRange.Value2 = $aArr
$bRange.Value2 = $bArr

# Apply the existing timestamp number format (style s="2") to all timestamp cells,
# including the newly added rows 94-97
$bRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"
